# Daily attendance processing - 2025-10-27 15:23:49
# Reorders the "Recorded By" email list (column G) for the affected session
# rows so that it matches the latest recording order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = "servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
    3  = "Amira.Sobhy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
    4  = "Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
    5  = "servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hananragab@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg"
    6  = "servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, gehanadel@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg"
    13 = "mariam.noureldin@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"
    14 = "Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"
    24 = "yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, marina_atef@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg"
    25 = "marina_atef@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg"
    26 = "gigi@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg"
    30 = "Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
    31 = "Amira.Sobhy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
    32 = "Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
    33 = "servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hananragab@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg"
    34 = "servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, gehanadel@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg"
    41 = "mariam.noureldin@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"
    42 = "Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"
    52 = "yasmin.m.senosy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, marina_atef@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg"
    53 = "marina_atef@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg"
    54 = "gigi@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg"
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
